$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 6).Value = 284 + ($r - 2)
}

$ws.Range("L12").Value = "stimuli/catch_25.jpg"
